# Updated symbol list (coinranking data refresh) - Sun Dec 18 15:50:05 UTC 2022
#
# This applies the cell-level updates to the "cryptos" worksheet: refreshed
# prices for many coins, and a handful of coins that changed rank/position
# in the table (so their Coin/Link/Price/Volume columns moved to a new row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    # Force text storage so numeric-looking strings (e.g. "246.30") are kept
    # exactly as text instead of being normalized into a Double, then restore
    # the default (unstyled) cell formatting used throughout the sheet.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell "D2" "246.30"
# Row 4
Set-TextCell "D4" "5.463"
# Row 5
Set-TextCell "D5" "0.05616"
# Row 6
Set-TextCell "D6" "6.470"
# Row 7
Set-TextCell "D7" "0.8046"
# Row 8
Set-TextCell "D8" "1.050"
# Row 9
Set-TextCell "D9" "0.1431"
# Row 10
Set-TextCell "D10" "0.07345"
# Row 11
Set-TextCell "D11" "0.03179"
# Row 12
Set-TextCell "D12" "0.02940"
# Row 13
Set-TextCell "D13" "0.09267"
# Row 14
Set-TextCell "D14" "0.001664"
# Row 15
Set-TextCell "D15" "3.218"
# Row 16
Set-TextCell "D16" "0.04740"
# Row 17
Set-TextCell "B17" "TigerCash"
Set-TextCell "C17" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell "D17" "0.006399"
Set-TextCell "E17" "16TigerCashTCH"
# Row 18
Set-TextCell "B18" "HotbitToken"
Set-TextCell "C18" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextCell "D18" "0.005065"
Set-TextCell "E18" "17HotbitTokenHTB"
# Row 19
Set-TextCell "D19" "0.001053"
# Row 20
Set-TextCell "B20" "NitroEx"
Set-TextCell "C20" "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextCell "D20" "0.0001504"
Set-TextCell "E20" "19NitroExNTX"
# Row 21
Set-TextCell "B21" "LEO"
Set-TextCell "C21" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell "D21" "3.985"
Set-TextCell "E21" "20LEOLEO"
# Row 22
Set-TextCell "B22" "GateToken"
Set-TextCell "C22" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextCell "D22" "3.385"
Set-TextCell "E22" "21GateTokenGT"
# Row 23
Set-TextCell "B23" "BTSEToken"
Set-TextCell "C23" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextCell "D23" "2.124"
Set-TextCell "E23" "22BTSETokenBTSE"
# Row 24
Set-TextCell "B24" "One"
Set-TextCell "C24" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell "D24" "0.01163"
Set-TextCell "E24" "23OneONEBestin24h"
# Row 26
Set-TextCell "D26" "0.1266"
Set-TextCell "E26" "25ProBitTokenPROB"
# Row 27
Set-TextCell "D27" "0.0002908"
# Row 40
Set-TextCell "D40" "0.04160"
# Row 41
Set-TextCell "D41" "0.006908"
# Row 42
Set-TextCell "B42" "CEJI"
Set-TextCell "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextCell "D42" "0.003509"
Set-TextCell "E42" "41CEJICEJI"
# Row 43
Set-TextCell "B43" "BKEXToken"
Set-TextCell "C43" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell "D43" "0.1037"
Set-TextCell "E43" "42BKEXTokenBKK"
# Row 44
Set-TextCell "D44" "0.008995"
# Row 45
Set-TextCell "D45" "0.00005663"
# Row 47
Set-TextCell "D47" "0.6818"
# Row 48
Set-TextCell "D48" "0.01735"

